$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Column B (accuracy values) per refreshed epoch run
$ws.Range("B2").Value = 0.40625
$ws.Range("B3").Value = 0.25
$ws.Range("B4").Value = 0.25
$ws.Range("B5").Value = 0.234375
$ws.Range("B6").Value = 0.21875
$ws.Range("B7").Value = 0.15625
$ws.Range("B8").Value = 0.1875
$ws.Range("B9").Value = 0.171875
$ws.Range("B10").Value = 0.1875
$ws.Range("B11").Value = 0.28125
$ws.Range("B12").Value = 0.21875
$ws.Range("B13").Value = 0.1875
$ws.Range("B14").Value = 0.1875
$ws.Range("B15").Value = 0.203125
$ws.Range("B16").Value = 0.15625
$ws.Range("B17").Value = 0.109375
$ws.Range("B18").Value = 0.109375
$ws.Range("B19").Value = 0.109375
$ws.Range("B20").Value = 0.109375
$ws.Range("B21").Value = 0.109375
$ws.Range("B22").Value = 0.109375
$ws.Range("B23").Value = 0.109375
$ws.Range("B24").Value = 0.109375
$ws.Range("B25").Value = 0.109375
$ws.Range("B26").Value = 0.109375
$ws.Range("B27").Value = 0.109375
$ws.Range("B28").Value = 0.109375
$ws.Range("B29").Value = 0.109375
$ws.Range("B30").Value = 0.109375
$ws.Range("B31").Value = 0.109375
$ws.Range("B32").Value = 0.125
$ws.Range("B33").Value = 0.125
$ws.Range("B34").Value = 0.125
$ws.Range("B35").Value = 0.125
$ws.Range("B43").Value = 0.140625
$ws.Range("B44").Value = 0.140625
$ws.Range("B45").Value = 0.140625
$ws.Range("B46").Value = 0.140625
$ws.Range("B47").Value = 0.140625
$ws.Range("B48").Value = 0.140625
$ws.Range("B49").Value = 0.140625
$ws.Range("B50").Value = 0.140625
$ws.Range("B51").Value = 0.140625
$ws.Range("B52").Value = 0.140625
$ws.Range("B53").Value = 0.140625
$ws.Range("B54").Value = 0.140625
$ws.Range("B55").Value = 0.140625
$ws.Range("B56").Value = 0.140625
$ws.Range("B57").Value = 0.140625
$ws.Range("B58").Value = 0.140625
$ws.Range("B59").Value = 0.140625
$ws.Range("B60").Value = 0.140625
$ws.Range("B61").Value = 0.140625
$ws.Range("B62").Value = 0.140625
$ws.Range("B63").Value = 0.140625
$ws.Range("B64").Value = 0.140625
$ws.Range("B65").Value = 0.140625
$ws.Range("B66").Value = 0.140625
$ws.Range("B67").Value = 0.140625
$ws.Range("B68").Value = 0.140625
$ws.Range("B69").Value = 0.140625
$ws.Range("B70").Value = 0.140625
$ws.Range("B71").Value = 0.140625
$ws.Range("B72").Value = 0.140625
$ws.Range("B73").Value = 0.140625
$ws.Range("B103").Value = 0.046875
$ws.Range("B104").Value = 0.046875
$ws.Range("B105").Value = 0.140625
$ws.Range("B106").Value = 0.125
$ws.Range("B107").Value = 0.078125
$ws.Range("B108").Value = 0.109375
$ws.Range("B109").Value = 0.0625
$ws.Range("B110").Value = 0.09375
$ws.Range("B112").Value = 0.125
$ws.Range("B113").Value = 0.078125
$ws.Range("B114").Value = 0.109375
$ws.Range("B117").Value = 0.046875
$ws.Range("B118").Value = 0.06557377049180328

# Column A rows 102-118 show the repr() of the DisplayOutputs object;
# the notebook was re-run so the object's memory address changed.
$newRepr = "<__main__.DisplayOutputs object at 0x7fd66042a250>"
for ($r = 102; $r -le 118; $r++) {
    $ws.Range("A$r").Value = $newRepr
}
